$d = $word.ActiveDocument

# Locate the sentence that needs to be split into two runs:
#   "Split up text into smaller units and count them again."
# becomes two runs:
#   "Split up text into smaller units and count them again " (trailing space)
#   "in parallel."
$findRange = $d.Content
$found = $findRange.Find.Execute(
    "Split up text into smaller units and count them again.",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    # Build a fresh Range over the matched span so InsertXML replaces its
    # contents (reusing the Find-affected range object inserts instead of
    # replacing).
    $target = $d.Range($findRange.Start, $findRange.End)

    $xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:body>' +
        '<w:p>' +
        '<w:r>' +
        '<w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr>' +
        '<w:t xml:space="preserve">Split up text into smaller units and count them again </w:t>' +
        '</w:r>' +
        '<w:r>' +
        '<w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr>' +
        '<w:t>in parallel.</w:t>' +
        '</w:r>' +
        '</w:p>' +
        '</w:body>' +
        '</w:document>' +
        '</pkg:xmlData></pkg:part></pkg:package>'

    [void]$target.InsertXML($xml)
}
